# Update gh-pages to output generated at 456a3b4
# This updates the "想去人数" (interest count) column F across all sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 528
$ws.Range("F3").Value = 10231
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 91
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 6867
$ws.Range("F8").Value = 649
$ws.Range("F9").Value = 0
$ws.Range("F12").Value = 12418
$ws.Range("F17").Value = 102
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 322
$ws.Range("F27").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 237
$ws.Range("F33").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 3621
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 117
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F47").Value = 274
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 4275
$ws.Range("F50").Value = 0

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 12
$ws.Range("F9").Value = 63
$ws.Range("F12").Value = 74
$ws.Range("F13").Value = 0
$ws.Range("F22").Value = 11
$ws.Range("F23").Value = 73
$ws.Range("F24").Value = 74
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("F29").Value = 0

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6363

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 10231
$ws.Range("F4").Value = 223
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 649
$ws.Range("F8").Value = 130
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 11880
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 33
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 191
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 1974
$ws.Range("F24").Value = 994
$ws.Range("F25").Value = 1468
$ws.Range("F26").Value = 0
$ws.Range("F28").Value = 2924
$ws.Range("F29").Value = 237
$ws.Range("F30").Value = 1975
$ws.Range("F31").Value = 105
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 1657
$ws.Range("F37").Value = 0
$ws.Range("F41").Value = 261
$ws.Range("F42").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 274
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 0
